$d = $word.ActiveDocument

# 1) Remove the stray "_GoBack" bookmark before "Controllare il funzionamento della data..."
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) Add a new "_GoBack" bookmark right after "(Es. Pag. 1/6)"
#    A zero-length range exactly on a paragraph's trailing boundary confuses
#    Bookmarks.Add in this host, so we briefly insert a one-character
#    placeholder, bookmark around it, then delete the placeholder again -
#    the bookmark collapses back to a zero-width mark in the right spot.
$r = $d.Content
$r.Find.Execute("(Es. Pag. 1/6)", $false, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Collapse(0)
    $r.InsertAfter("X")
    $d.Bookmarks.Add("_GoBack", $r)
    $placeholder = $d.Range($r.Start, $r.End)
    $placeholder.Text = ""
}

# 3) Highlight the "Fornitore" field paragraph yellow instead of cyan
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Aggiungere campo (stringa 50 caratteri) nominato*") {
        $p.Range.Font.HighlightColorIndex = 7
    }
}
